# starter region 4 results
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("model_rerun_results_tracker")

# Rows 322-329: still "not run" -> placeholders "NA" for H/I/J/L, "T" for K
$naRows = 322..329
foreach ($r in $naRows) {
    $ws.Cells.Item($r, 8).Value = "NA"   # H num_vehicles
    $ws.Cells.Item($r, 9).Value = "NA"   # I cumu_time
    $ws.Cells.Item($r, 10).Value = "NA"  # J cumu_dist
    $ws.Cells.Item($r, 11).Value = "T"   # K solution_file
    $ws.Cells.Item($r, 12).Value = "NA"  # L map_file
}

# Rows 330-353: starter region 4 results filled in
$results = @{
    330 = @(7, 2329, 1126.86)
    331 = @(8, 2695, 1519.24)
    332 = @(8, 2912, 1558.6)
    333 = @(8, 2658, 1465.3)
    334 = @(6, 2537, 1305.26)
    335 = @(6, 2674, 1343.06)
    336 = @(7, 2631, 1367.83)
    337 = @(6, 2726, 1396.91)
    338 = @(5, 2327, 1104.15)
    339 = @(7, 2594, 1350.12)
    340 = @(7, 2534, 1304.95)
    341 = @(7, 2736, 1412.54)
    342 = @(4, 2135, 985.32)
    343 = @(5, 2325, 1127.87)
    344 = @(6, 2561, 1283.95)
    345 = @(5, 2298, 1102.31)
    346 = @(5, 2327, 1104.15)
    347 = @(7, 2582, 1341.42)
    348 = @(7, 2534, 1304.95)
    349 = @(7, 2736, 1412.54)
    350 = @(4, 2135, 985.32)
    351 = @(5, 2325, 1127.87)
    352 = @(6, 2561, 1283.95)
    353 = @(5, 2360, 1138.05)
}

foreach ($r in 330..353) {
    $vals = $results[$r]
    $ws.Cells.Item($r, 8).Value = $vals[0]   # H num_vehicles
    $ws.Cells.Item($r, 9).Value = $vals[1]   # I cumu_time
    $ws.Cells.Item($r, 10).Value = $vals[2]  # J cumu_dist
    $ws.Cells.Item($r, 11).Value = "T"       # K solution_file
    $ws.Cells.Item($r, 12).Value = "T"       # L map_file
}

# Update the view state to match where the user scrolled/selected after
# entering this batch of results.
$ws.Application.ActiveWindow.ScrollRow = 339
$ws.Application.ActiveWindow.ScrollColumn = 3
$ws.Range("H354").Select()
